$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with the new values (34 columns A:AH)
$data = New-Object 'object[,]' 4,34
$data[0,0] = 45172.50694444445
$data[0,1] = 21.139
$data[0,2] = 14.405
$data[0,3] = 4.093
$data[0,4] = 44.473
$data[0,5] = 36.592
$data[0,6] = 16.635
$data[0,7] = 54.157
$data[0,8] = 25.596
$data[0,9] = 10.793
$data[0,10] = 16.568
$data[0,11] = 17.656
$data[0,12] = 18.451
$data[0,13] = 5.311
$data[0,14] = 16.542
$data[0,15] = 23.194
$data[0,16] = 13.966
$data[0,17] = 3.726
$data[0,18] = 2.455
$data[0,19] = 244.445
$data[0,20] = 46.035
$data[0,21] = 15.269
$data[0,22] = 30.403
$data[0,23] = 15.788
$data[0,24] = 2.61
$data[0,25] = 26.878
$data[0,26] = 13.487
$data[0,27] = 12.177
$data[0,28] = 14.251
$data[0,29] = 18.232
$data[0,30] = 3.64
$data[0,31] = 47.882
$data[0,32] = 8.427
$data[0,33] = 19.09
$data[1,0] = 45172.51388888889
$data[1,1] = 3.363
$data[1,2] = 1.826
$data[1,3] = 1.371
$data[1,4] = 6.725
$data[1,5] = 5.458
$data[1,6] = 2.648
$data[1,7] = 16.197
$data[1,8] = 4.072
$data[1,9] = 1.589
$data[1,10] = 2.313
$data[1,11] = 2.755
$data[1,12] = 2.718
$data[1,13] = 0.863
$data[1,14] = 2.632
$data[1,15] = 3.642
$data[1,16] = 2.544
$data[1,17] = 1.498
$data[1,18] = 0.6899999999999999
$data[1,19] = 32.822
$data[1,20] = 7.777
$data[1,21] = 2.429
$data[1,22] = 4.876
$data[1,23] = 2.58
$data[1,24] = 0.671
$data[1,25] = 7.098
$data[1,26] = 2.146
$data[1,27] = 2.161
$data[1,28] = 2.485
$data[1,29] = 2.755
$data[1,30] = 1.294
$data[1,31] = 15.157
$data[1,32] = 1.187
$data[1,33] = 3.042
$data[2,0] = 45172.52083333334
$data[2,1] = 20.178
$data[2,2] = 14.701
$data[2,3] = 1.447
$data[2,4] = 43.528
$data[2,5] = 35.895
$data[2,6] = 15.879
$data[2,7] = 57.341
$data[2,8] = 24.432
$data[2,9] = 10.748
$data[2,10] = 16.07
$data[2,11] = 17.551
$data[2,12] = 18.407
$data[2,13] = 5.07
$data[2,14] = 15.79
$data[2,15] = 22.383
$data[2,16] = 13.38
$data[2,17] = 1.107
$data[2,18] = 0.93
$data[2,19] = 233.011
$data[2,20] = 44.036
$data[2,21] = 14.575
$data[2,22] = 29.479
$data[2,23] = 15.659
$data[2,24] = 2.33
$data[2,25] = 28.177
$data[2,26] = 12.874
$data[2,27] = 11.489
$data[2,28] = 13.488
$data[2,29] = 18.36
$data[2,30] = 0.784
$data[2,31] = 51.705
$data[2,32] = 8.161
$data[2,33] = 18.222
$data[3,0] = 45172.52777777778
$data[3,1] = 12.97
$data[3,2] = 9.4
$data[3,3] = 0.99
$data[3,4] = 27.96
$data[3,5] = 23.02
$data[3,6] = 10.21
$data[3,7] = 41.59
$data[3,8] = 15.71
$data[3,9] = 6.91
$data[3,10] = 10.27
$data[3,11] = 11.3
$data[3,12] = 11.82
$data[3,13] = 3.26
$data[3,14] = 10.15
$data[3,15] = 14.4
$data[3,16] = 8.66
$data[3,17] = 0.8
$data[3,18] = 0.61
$data[3,19] = 147.19
$data[3,20] = 28.46
$data[3,21] = 9.369999999999999
$data[3,22] = 19
$data[3,23] = 10.1
$data[3,24] = 1.55
$data[3,25] = 19.78
$data[3,26] = 8.279999999999999
$data[3,27] = 7.41
$data[3,28] = 8.699999999999999
$data[3,29] = 11.81
$data[3,30] = 0.5600000000000001
$data[3,31] = 37.84
$data[3,32] = 5.22
$data[3,33] = 11.71
$ws.Range("A2:AH5").Value = $data

# Adjust specific column widths from 7 -> 8 characters (stored width 7 -> 8)
# ColumnWidth property adds ~0.8333 padding internally; 7.166666666666667 -> stored width 8
$ws.Range("J:J").ColumnWidth = 7.166666666666667
$ws.Range("Q:Q").ColumnWidth = 7.166666666666667
$ws.Range("AA:AC").ColumnWidth = 7.166666666666667

# Remove row 6 (data reduced from 5 data rows to 4 data rows); dimension auto-updates to A1:AH5
$ws.Rows.Item(6).Delete()

Write-Host "edit complete"
